# Update correlation output data (Pearson table) for
# pearson_tables/t2m_netherlands-2-3.xlsx
# Columns: B = SST, C = MSLP, D = Z500
# Rows 2-9 correspond to SCA_N, EA_N, ENSO-mei_N, NAO_N, SCA_P, EA_P, ENSO-mei_P, NAO_P

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = -0.598108242147238
$ws.Range("C2").Value = 0.6534443712313868
$ws.Range("D2").Value = -0.6789664725740319

$ws.Range("B3").Value = -0.7469354083221468
$ws.Range("C3").Value = 0.700834921232218
$ws.Range("D3").Value = -0.8017516749261175

$ws.Range("B4").Value = -0.6605507161631865
$ws.Range("C4").Value = -0.6630486195838414
$ws.Range("D4").Value = -0.6803745916549853

$ws.Range("B5").Value = -0.8448773721439531
$ws.Range("C5").Value = 0.6056541206630182
$ws.Range("D5").Value = 0.674227746103139

$ws.Range("B6").Value = 0.7470809345076589
$ws.Range("C6").Value = 0.6707235140136254
$ws.Range("D6").Value = -0.6229360296395556

$ws.Range("B7").Value = -0.7162815057609989
$ws.Range("C7").Value = -0.6555241559669299
$ws.Range("D7").Value = 0.7713944330297057

$ws.Range("B8").Value = -0.7873339243053709
$ws.Range("C8").Value = 0.4936852292443559
$ws.Range("D8").Value = 0.6384932586089653

$ws.Range("B9").Value = -0.7423493649827849
$ws.Range("C9").Value = -0.6041343640987742
$ws.Range("D9").Value = 0.6510576106228616
